$d = $word.ActiveDocument

function CleanText($range) {
    return $range.Text.TrimEnd([char]13, [char]7)
}

# PowerShell's -eq / -ceq / -clike all fold case in this host, so match
# case-sensitively via a case-sensitive regex instead (needed to tell
# "Cadastrar Paciente" apart from the "CADASTRAR PACIENTE" heading, etc.).
function TextEquals($a, $b) {
    $pattern = "(?-i)^" + [regex]::Escape($b) + "$"
    return $a -cmatch $pattern
}

# --- 1) Add "Cadastrar Prontuário" bullet right after "Cadastrar Consulta"
#        in the top "Necessidades" list (numId=1, ilvl=0 bullet list that
#        also holds "Cadastrar Paciente" / "Cadastrar Médico"). ---
$found1 = $false
for ($i = 1; $i -le $d.Paragraphs.Count -and -not $found1; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = CleanText $p.Range
    if ((TextEquals $t "Cadastrar Consulta") -and $p.Range.ListFormat.ListLevelNumber -eq 1) {
        $prev = $d.Paragraphs.Item($i - 1)
        if (TextEquals (CleanText $prev.Range) "Cadastrar Médico") {
            $p.Range.InsertParagraphAfter()
            $newPar = $d.Paragraphs.Item($i + 1)
            $newPar.Range.Text = "Cadastrar Prontuário"
            $found1 = $true
        }
    }
}

# --- 2) Add a second empty "no-numbering, 720-indent" spacer paragraph
#        right after the existing one (between the "Necessidades" list and
#        the blank line that precedes "CADASTRAR PACIENTE"). ---
$found2 = $false
for ($i = 1; $i -le $d.Paragraphs.Count -and -not $found2; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = CleanText $p.Range
    if ((TextEquals $t "") -and $p.Range.ParagraphFormat.LeftIndent -eq 36) {
        $p.Range.InsertParagraphAfter()
        $found2 = $true
    }
}

# --- 3) Add "Cadastrar contato, telefones e e-mails" sub-bullet right
#        after the "Cadastrar Paciente" bullet that opens the
#        "CADASTRAR PACIENTE" section's Checklist (numId=3, ilvl=0 -> new
#        paragraph at ilvl=1). ---
$found3 = $false
$seenHeading = $false
for ($i = 1; $i -le $d.Paragraphs.Count -and -not $found3; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = CleanText $p.Range
    if (-not $seenHeading) {
        if (TextEquals $t "CADASTRAR PACIENTE") {
            $seenHeading = $true
        }
        continue
    }
    if ((TextEquals $t "Cadastrar Paciente") -and $p.Range.ListFormat.ListLevelNumber -eq 1) {
        $p.Range.InsertParagraphAfter()
        $newPar = $d.Paragraphs.Item($i + 1)
        $newPar.Range.Text = "Cadastrar contato, telefones e e-mails"
        $newPar.Range.Font.Color = 16711775
        $newPar.Range.Font.Size = 9
        $newPar.Range.Font.SizeBi = 9
        $newPar.Range.ListFormat.ListLevelNumber = 2
        $found3 = $true
    }
}
